$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = [double]"25.01000000000047"
$ws.Range("H2").Value = [double]"0.0003982719897474007"
$ws.Range("I2").Value = [double]"0.0003982719897474007"
$ws.Range("L2").Value = [double]"42.98414391871968"
$ws.Range("M2").Value = "[17.378484808620854, 68.58980302881851]"
$ws.Range("N2").Value = [double]"0.001502037087650843"
$ws.Range("O2").Value = [double]"0.001502037087650843"
$ws.Range("P2").Value = [double]"1.86797401024258"
$ws.Range("Q2").Value = "[1.1761317842268104, 2.5598162362583503]"
$ws.Range("R2").Value = [double]"2.112949634813788e-06"
$ws.Range("S2").Value = [double]"2.112949634813788e-06"
$ws.Range("T2").Value = [double]"61.13419945698831"
$ws.Range("U2").Value = "[46.67680783070293, 75.59159108327368]"
$ws.Range("V2").Value = [double]"6.174816213899703e-11"
$ws.Range("W2").Value = [double]"6.174816213899703e-11"
$ws.Range("X2").Value = [double]"17.57459459459492"
$ws.Range("Y2").Value = [double]"14.82074074074102"
$ws.Range("Z2").Value = [double]"20.32844844844882"
$ws.Range("F3").Value = [double]"25.01000000000047"
$ws.Range("H3").Value = [double]"0.000156419883225567"
$ws.Range("I3").Value = [double]"0.000156419883225567"
$ws.Range("L3").Value = [double]"37.57249527459734"
$ws.Range("M3").Value = "[16.824218244908586, 58.3207723042861]"
$ws.Range("N3").Value = [double]"0.0006854150987423147"
$ws.Range("O3").Value = [double]"0.0006854150987423147"
$ws.Range("P3").Value = [double]"1.276763380738194"
$ws.Range("Q3").Value = "[0.6603948521059619, 1.8931319093704255]"
$ws.Range("R3").Value = [double]"0.0001359863332710542"
$ws.Range("S3").Value = [double]"0.0001359863332710542"
$ws.Range("T3").Value = [double]"50.62503847591687"
$ws.Range("U3").Value = "[38.92788254658856, 62.322194405245185]"
$ws.Range("V3").Value = [double]"3.19406723292559e-11"
$ws.Range("W3").Value = [double]"3.19406723292559e-11"
$ws.Range("X3").Value = [double]"19.92788788788826"
$ws.Range("Y3").Value = [double]"17.47445445445478"
$ws.Range("Z3").Value = [double]"22.38132132132175"
$ws.Range("B4").Value = [double]"0"
$ws.Range("F4").Value = [double]"25.01000000000047"
$ws.Range("H4").Value = [double]"0.01812264687271825"
$ws.Range("I4").Value = [double]"0.01812264687271825"
$ws.Range("L4").Value = [double]"26.54453315541483"
$ws.Range("M4").Value = "[2.7471324515652285, 50.34193385926444]"
$ws.Range("N4").Value = [double]"0.02961774843110376"
$ws.Range("O4").Value = [double]"0.02961774843110376"
$ws.Range("P4").Value = [double]"1.188710733790733"
$ws.Range("Q4").Value = "[-0.03144737390980801, 2.408868841491273]"
$ws.Range("R4").Value = [double]"0.05594192536274689"
$ws.Range("S4").Value = [double]"0.05594192536274689"
$ws.Range("T4").Value = [double]"49.05709082126936"
$ws.Range("U4").Value = "[36.17446938537806, 61.93971225716065]"
$ws.Range("V4").Value = [double]"1.047149922328572e-09"
$ws.Range("W4").Value = [double]"1.047149922328572e-09"
$ws.Range("X4").Value = [double]"20.27837837837876"
$ws.Range("Y4").Value = [double]"15.42158158158187"
$ws.Range("Z4").Value = [double]"25.13517517517565"
$ws.Range("F5").Value = [double]"25.01000000000047"
$ws.Range("H5").Value = [double]"7.825228864888345e-08"
$ws.Range("I5").Value = [double]"7.825228864888345e-08"
$ws.Range("L5").Value = [double]"62.49062592073967"
$ws.Range("M5").Value = "[37.68070178278542, 87.30055005869393]"
$ws.Range("N5").Value = [double]"7.209520287654314e-06"
$ws.Range("O5").Value = [double]"7.209520287654314e-06"
$ws.Range("P5").Value = [double]"0.8616580451287312"
$ws.Range("Q5").Value = "[0.47171060864711656, 1.251605481610346]"
$ws.Range("R5").Value = [double]"5.58604139038188e-05"
$ws.Range("S5").Value = [double]"5.58604139038188e-05"
$ws.Range("T5").Value = [double]"66.93656880092415"
$ws.Range("U5").Value = "[54.16942173639957, 79.70371586544873]"
$ws.Range("V5").Value = [double]"9.192646643896296e-14"
$ws.Range("W5").Value = [double]"9.192646643896296e-14"
$ws.Range("X5").Value = [double]"21.58020020020061"
$ws.Range("Y5").Value = [double]"20.02802802802841"
$ws.Range("Z5").Value = [double]"23.1323723723728"
$ws.Range("F6").Value = [double]"25.01000000000047"
$ws.Range("H6").Value = [double]"8.621327084190611e-05"
$ws.Range("I6").Value = [double]"8.621327084190611e-05"
$ws.Range("L6").Value = [double]"51.39132260500694"
$ws.Range("M6").Value = "[22.10377857891021, 80.67886663110367]"
$ws.Range("N6").Value = [double]"0.0009597114593944411"
$ws.Range("O6").Value = [double]"0.0009597114593944411"
$ws.Range("P6").Value = [double]"0.6855527512338089"
$ws.Range("Q6").Value = "[0.16981581911296306, 1.2012896833546547]"
$ws.Range("R6").Value = [double]"0.01032121786584717"
$ws.Range("S6").Value = [double]"0.01032121786584717"
$ws.Range("T6").Value = [double]"64.44339563248705"
$ws.Range("U6").Value = "[49.42599112998333, 79.46080013499078]"
$ws.Range("V6").Value = [double]"4.073075210442312e-11"
$ws.Range("W6").Value = [double]"4.073075210442312e-11"
$ws.Range("X6").Value = [double]"22.2811811811816"
$ws.Range("Y6").Value = [double]"20.22830830830869"
$ws.Range("Z6").Value = [double]"24.33405405405451"
$ws.Range("F7").Value = [double]"25.01000000000047"
$ws.Range("H7").Value = [double]"0.0002285987757412178"
$ws.Range("I7").Value = [double]"0.0002285987757412178"
$ws.Range("L7").Value = [double]"40.58880416474915"
$ws.Range("M7").Value = "[16.692169772084156, 64.48543855741414]"
$ws.Range("N7").Value = [double]"0.001337666659656733"
$ws.Range("O7").Value = [double]"0.001337666659656733"
$ws.Range("P7").Value = [double]"0.4717106086471166"
$ws.Range("Q7").Value = "[-0.15723686954903737, 1.1006580868432705]"
$ws.Range("R7").Value = [double]"0.1378868762615149"
$ws.Range("S7").Value = [double]"0.1378868762615149"
$ws.Range("T7").Value = [double]"58.65764498973205"
$ws.Range("U7").Value = "[46.05966226721188, 71.25562771225222]"
$ws.Range("V7").Value = [double]"3.738565013122752e-12"
$ws.Range("W7").Value = [double]"3.738565013122752e-12"
$ws.Range("X7").Value = [double]"23.13237237237281"
$ws.Range("Y7").Value = [double]"20.62886886886926"
$ws.Range("Z7").Value = [double]"25.63587587587635"
$ws.Range("F8").Value = [double]"23.70000000000027"
$ws.Range("H8").Value = [double]"0.0004369418129007707"
$ws.Range("I8").Value = [double]"0.0004369418129007707"
$ws.Range("L8").Value = [double]"44.98597484602176"
$ws.Range("M8").Value = "[18.67221184951947, 71.29973784252405]"
$ws.Range("N8").Value = [double]"0.001253372505753525"
$ws.Range("O8").Value = [double]"0.001253372505753525"
$ws.Range("P8").Value = [double]"0.1823947686768852"
$ws.Range("Q8").Value = "[-0.4968685077749617, 0.8616580451287321]"
$ws.Range("R8").Value = [double]"0.5912960703101664"
$ws.Range("S8").Value = [double]"0.5912960703101664"
$ws.Range("T8").Value = [double]"60.83966978219445"
$ws.Range("U8").Value = "[45.84281717568281, 75.83652238870609]"
$ws.Range("V8").Value = [double]"1.946856009737985e-10"
$ws.Range("W8").Value = [double]"1.946856009737985e-10"
$ws.Range("X8").Value = [double]"23.01201201201227"
$ws.Range("Y8").Value = [double]"20.44984984985008"
$ws.Range("Z8").Value = [double]"25.57417417417446"
$ws.Range("B9").Value = [double]"0"
$ws.Range("F9").Value = [double]"23.70000000000027"
$ws.Range("H9").Value = [double]"0.001527498002526939"
$ws.Range("I9").Value = [double]"0.001527498002526939"
$ws.Range("L9").Value = [double]"40.90706249127501"
$ws.Range("M9").Value = "[14.691421744657546, 67.12270323789248]"
$ws.Range("N9").Value = [double]"0.002958846799122838"
$ws.Range("O9").Value = [double]"0.002958846799122838"
$ws.Range("P9").Value = [double]"0.3207632138800394"
$ws.Range("Q9").Value = "[-0.4339737599553457, 1.0755001877154244]"
$ws.Range("R9").Value = [double]"0.3965363059738554"
$ws.Range("S9").Value = [double]"0.3965363059738554"
$ws.Range("T9").Value = [double]"57.96824028242474"
$ws.Range("U9").Value = "[42.869893892860176, 73.0665866719893]"
$ws.Range("V9").Value = [double]"8.460754319372654e-10"
$ws.Range("W9").Value = [double]"8.460754319372654e-10"
$ws.Range("X9").Value = [double]"22.49009009009034"
$ws.Range("Y9").Value = [double]"19.64324324324346"
$ws.Range("Z9").Value = [double]"25.33693693693722"
$ws.Range("F10").Value = [double]"23.70000000000027"
$ws.Range("H10").Value = [double]"0.0003808817154613564"
$ws.Range("I10").Value = [double]"0.0003808817154613564"
$ws.Range("L10").Value = [double]"41.97478785762807"
$ws.Range("M10").Value = "[17.25926021951777, 66.69031549573836]"
$ws.Range("N10").Value = [double]"0.001339228719998253"
$ws.Range("O10").Value = [double]"0.001339228719998253"
$ws.Range("P10").Value = [double]"0.3333421634439624"
$ws.Range("Q10").Value = "[-0.3207632138800385, 0.9874475407679633]"
$ws.Range("R10").Value = [double]"0.3101805318316155"
$ws.Range("S10").Value = [double]"0.3101805318316155"
$ws.Range("T10").Value = [double]"63.41598053089972"
$ws.Range("U10").Value = "[49.600403166494836, 77.2315578953046]"
$ws.Range("V10").Value = [double]"5.730305119300283e-12"
$ws.Range("W10").Value = [double]"5.730305119300283e-12"
$ws.Range("X10").Value = [double]"22.44264264264289"
$ws.Range("Y10").Value = [double]"19.9753753753756"
$ws.Range("Z10").Value = [double]"24.90990990991019"
